$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Single-cell corrections in existing rows ---
# E8: "4" -> "2B"
$ws.Range("E8").Value = "2B"

# I39: "4" -> "2B", L39: "4" -> "3"
$ws.Range("I39").Value = "2B"

# L39 becomes the numeric-looking text "3"; force text storage so it keeps
# referencing a shared string instead of becoming a literal number.
$ws.Range("L39").NumberFormat = "@"
$ws.Range("L39").Value = "3"
$ws.Range("L39").ClearFormats()

# --- Insert a new row at 221, pushing existing rows 221:251 down to 222:252 ---
$ws.Rows(221).Insert()

# --- Populate the newly inserted row 221 ---
$ws.Range("A221").Value = 221
$ws.Range("C221").Value = "1B"
$ws.Range("E221").Value = "2B"
$ws.Range("F221").Value = "2B"
$ws.Range("G221").Value = "2A"
$ws.Range("H221").Value = "1B"
$ws.Range("P221").Value = "2A"
$ws.Range("Q221").Value = "2A"
$ws.Range("R221").Value = "1B"

# These values ("4" / "3") look like plain numbers, so they must be forced
# to text storage (shared-string cells), matching the rest of the sheet.
$textCells = @("B221", "D221", "I221", "J221", "K221", "L221", "M221", "N221", "O221")
$textValues = @{
    "B221" = "4"
    "D221" = "3"
    "I221" = "4"
    "J221" = "4"
    "K221" = "3"
    "L221" = "4"
    "M221" = "4"
    "N221" = "4"
    "O221" = "4"
}
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $textValues[$addr]
    $ws.Range($addr).ClearFormats()
}
